# Apply marksheet corrections: update correct/total marks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# B11: Marking (right) weight 3 -> 5
$ws.Range("B11").Value = 5

# B12: Total (right) marks 30 -> 50
$ws.Range("B12").Value = 50

# E12: Corr/total marks text "18/84" -> "50/140"
$ws.Range("E12").Value = "50/140"
